{"js": "// The closing \"Eu, ... RG: {{ num_rg }}\" block keeps its wording; only the\n// signature line's date placeholder changes from {{ data }} to\n// {{ data_assinatura }} (e.g. \"{{ cidade }}, {{ data }}.\" -> \"{{ cidade }}, {{ data_assinatura }}.\").\nconst body = context.document.body;\nconst results = body.search(\"{{ data }}\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"{{ data_assinatura }}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The closing \"Eu, ... RG: {{ num_rg }}\" block keeps its wording; only the\n# signature line's date placeholder changes from {{ data }} to\n# {{ data_assinatura }} (e.g. \"{{ cidade }}, {{ data }}.\" -> \"{{ cidade }}, {{ data_assinatura }}.\").\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"{{ data }}\"\n$find.Replacement.Text = \"{{ data_assinatura }}\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $null, $null, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2) | Out-Null\n"}
